# Generate Report for Handback
# Update the timestamp cells recorded on the "Overview", "zh-cn" and "de-de"
# sheets to reflect the freshly generated handback report.

$wb = $excel.ActiveWorkbook

# -- Overview sheet: "Latest HO Xliff Generate Date" (column G, row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 11:01:45"

# -- zh-cn sheet: "Correspond Handoff Datetime" (H2) and
#                 "Correspond Handback DateTime" (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 11:01:39"
$wsZhCn.Range("K2").Value = "2016-08-17 11:01:56"

# -- de-de sheet: "Correspond Handoff Datetime" (H2) and
#                 "Correspond Handback DateTime" (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-17 11:01:45"
$wsDeDe.Range("K2").Value = "2016-08-17 11:02:10"
